$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Not living with children.enfants"
$ws.Range("C1").Value = "Living with children.enfants"
$ws.Range("D1").Value = "Not known / missing.enfants"
$ws.Range("E1").Value = "Total.enfants"
